# Update column G ("K") values for rows 2-27 on Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$values = @{
    2  = 1
    3  = 2
    4  = 3
    5  = 3
    6  = 4
    7  = 5
    8  = 4
    9  = 5
    10 = 4
    11 = 2
    12 = 2
    13 = 3
    14 = 1
    15 = 0
    16 = 2
    17 = 2
    18 = 1
    19 = 0
    20 = 0
    21 = 0
    22 = 1
    23 = 0
    24 = 0
    25 = 0
    26 = 1
    27 = 1
}

foreach ($row in $values.Keys) {
    $ws.Range("G$row").Value = $values[$row]
}
